# Apply the "P2 Industry rates" update:
#  - Increase every "Expansion Hiring %" (column G, rows 2-23) value by 0.5
#  - Move the active cell selection on Sheet1 to A24

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G
    $cell.Value2 = $cell.Value2 + 0.5
}

# Update the selected cell to match the saved view state
$ws.Range("A24").Select()

$wb.Save()
